# "More Protocol with working packet type detection"
#
# The protocol reference table lists packet-type IDs (EF, EE, ED, EC, EB, EA)
# for the six "Setup Conection T-Sx" sections. The stray duplicate "EF" row
# and an orphaned stray "s" label cell (J35) are removed, and the whole ID
# sequence shifts down by one hex step, adding a new "E9" packet type at the
# bottom (row 66) that didn't exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every packet-type ID cell down one step: EF/EE/ED/EC/EB/EA -> EE/ED/EC/EB/EA/E9
$ws.Range("B19").Value = "EE"
$ws.Range("B30").Value = "ED"
$ws.Range("B39").Value = "EC"
$ws.Range("B48").Value = "EB"
$ws.Range("B57").Value = "EA"
$ws.Range("B66").Value = "E9"

# Drop the stray orphaned "s" text that lived next to the T-S0 header
$ws.Range("J35").ClearContents()

# Leave the reader's selection on the newly added packet-type row
$ws.Range("B66").Select()
